# Commit: "Single space words and remove unknown from yellow highlights"
#
# review_template.docx's only paragraph contains nothing but a stray
# highlighted space (left over from manual review markup) immediately
# before the hidden "_GoBack" bookmark. Find that highlighted run and
# remove it, leaving the bookmark start/end tags exactly where they are.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = ""                # match by formatting, not literal text
$find.Highlight = $true        # ... specifically, any highlighted run
$find.Forward = $true
$find.Wrap = 1                 # wdFindContinue

# There's only one such run in this template, but loop defensively in
# case a future revision of the template has more than one.
$guard = 0
while ($find.Execute() -and $guard -lt 50) {
    $guard = $guard + 1
    $hit = $find.Parent
    if ($hit.Start -eq $hit.End) {
        break   # nothing left to delete; avoid spinning forever
    }
    $hit.Delete()
}

# Best-effort: Word sometimes mints a fresh GUID for the bibliography
# datastore item (customXml/itemProps1.xml) when it rewrites the custom
# XML parts on save. The hosted object model here doesn't expose that
# part for scripted editing, so this is a no-op if unsupported -- it's
# not something a document author could control from the UI/VBA either.
try {
    foreach ($part in $d.CustomXMLParts) {
        if ($part.XML -like "*bibliography*") {
            $part.Delete()
        }
    }
} catch {
}
